$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns L:N, copying the header style (bold, border, centered) from K1
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)

$ws.Range("L1").Value = "Contrato"
$ws.Range("M1").Value = "Registro"
$ws.Range("N1").Value = "Modalidade"

# Protocolo (E) and CPF (G) columns hold numeric-looking strings that must stay text
# (to preserve leading zeros, e.g. CPF "04466559457"). Mark as text before writing,
# then restore the default "Normal" style so no stray formatting is left behind.
$ws.Range("E2:E3").NumberFormat = "@"
$ws.Range("G2:G3").NumberFormat = "@"

# Row 2 values
$ws.Cells.Item(2, 1).Value = "422380 - YOU ASSISTÊNCIA MÉDICA LTDA"
$ws.Cells.Item(2, 2).Value = "18-04-2023"
$ws.Cells.Item(2, 3).Value = 45022.40521990741
$ws.Cells.Item(2, 4).Value = 12159014
$ws.Cells.Item(2, 5).Value = "8595238"
$ws.Cells.Item(2, 6).Value = "CARLOS RODRIGO CHAGAS GITIRANA"
$ws.Cells.Item(2, 7).Value = "04466559457"
$ws.Cells.Item(2, 8).Value = "No dia 22/03/2023 eu solicitei autorização para o exame de colonoscopia com biópsia, já se passaram 10 dias úteis e ainda não foi autorizado "
$ws.Cells.Item(2, 9).Value = 4
$ws.Cells.Item(2, 10).Value = "NO"
$ws.Cells.Item(2, 11).Value = "Assistencial"

# Row 3 values
$ws.Cells.Item(3, 1).Value = "422380 - YOU ASSISTÊNCIA MÉDICA LTDA"
$ws.Cells.Item(3, 2).Value = "18-04-2023"
$ws.Cells.Item(3, 3).Value = 45022.71378472223
$ws.Cells.Item(3, 4).Value = 12160376
$ws.Cells.Item(3, 5).Value = "8596884"
$ws.Cells.Item(3, 6).Value = "CRISTIANO DE OLIVEIRA CARNEIRO"
$ws.Cells.Item(3, 7).Value = "87895145487"
$ws.Cells.Item(3, 8).Value = "Fiz a migração do plano Hapvida para o you saúde,  e foi nos garantido a redução das carências pela corretora de saúde e quando precisei de exames não foi autorizada. Entrei em contato com a you saúde e a administradora Sindfort e não deram solução.`r`nObs entreguei todos os documentos solicitados no início de janeiro e nada. Estou sendo enganado e lesado ."
$ws.Cells.Item(3, 9).Value = 4
$ws.Cells.Item(3, 10).Value = "NO"
$ws.Cells.Item(3, 11).Value = "Assistencial"

# Restore default formatting for the Protocolo/CPF cells (keeps them text, drops the
# temporary "@" number format so the style matches the rest of the untouched columns).
$ws.Range("E2:E3").Style = "Normal"
$ws.Range("G2:G3").Style = "Normal"

# The multi-line description in H3 triggers an automatic row-height bump; reset and
# autofit so row 3 keeps its default (no explicit height override), like the source.
$ws.Rows(3).RowHeight = 15
$ws.Rows(3).EntireRow.AutoFit()
